$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.5216745459112635
$ws.Range("B2").Value = 59.99999999999999
$ws.Range("A3").Value = -0.5096733457912509
$ws.Range("B3").Value = 60.00000000000002
$ws.Range("A4").Value = -0.4976721456712409
$ws.Range("B4").Value = 60.00000000000002
$ws.Range("A5").Value = -0.4856709455512273
$ws.Range("B5").Value = 60.00000000000002
$ws.Range("A6").Value = -0.4736697454312156
$ws.Range("B6").Value = 59.99999999999999
$ws.Range("A7").Value = -0.4616685453112032
$ws.Range("B7").Value = 59.99999999999999
$ws.Range("A8").Value = -0.4496673451911906
$ws.Range("B8").Value = 60.00000000000002
$ws.Range("A9").Value = -0.4376661450711794
$ws.Range("B9").Value = 59.99999999999999
$ws.Range("A10").Value = -0.4256649449511669
$ws.Range("B10").Value = 60.00000000000002
$ws.Range("A11").Value = -0.4136637448311555
$ws.Range("B11").Value = 59.99999999999999
$ws.Range("A12").Value = -0.4016625447111435
$ws.Range("B12").Value = 60.00000000000002
$ws.Range("A13").Value = -0.3896613445911514
$ws.Range("B13").Value = 59.99999999999999
$ws.Range("A14").Value = -0.3776601444711178
$ws.Range("B14").Value = 60.00000000000002
$ws.Range("A15").Value = -0.3656589443511068
$ws.Range("B15").Value = 60.00000000000002
$ws.Range("A16").Value = -0.3536577442310964
$ws.Range("B16").Value = 60.00000000000002
$ws.Range("A17").Value = -0.3416565441110834
$ws.Range("B17").Value = 59.99999999999999
$ws.Range("A18").Value = -0.3296553439910707
$ws.Range("B18").Value = 60.00000000000002
$ws.Range("A19").Value = -0.31765414387106
$ws.Range("B19").Value = 60.00000000000002
$ws.Range("A20").Value = -0.3056529437510471
$ws.Range("B20").Value = 59.99999999999999
$ws.Range("A21").Value = -0.2936517436310367
$ws.Range("B21").Value = 60.00000000000002
$ws.Range("A22").Value = -0.2816505435110241
$ws.Range("B22").Value = 60.00000000000002
$ws.Range("A23").Value = -0.2696493433910113
$ws.Range("B23").Value = 59.99999999999999
$ws.Range("A24").Value = -0.2576481432710014
$ws.Range("B24").Value = 60.00000000000002
$ws.Range("A25").Value = -0.2456469431509883
$ws.Range("B25").Value = 59.99999999999999
$ws.Range("A26").Value = -0.2336457430309758
$ws.Range("B26").Value = 60.00000000000002
$ws.Range("A27").Value = -0.2216445429109649
$ws.Range("B27").Value = 59.99999999999999
$ws.Range("A28").Value = -0.2096433427909529
$ws.Range("B28").Value = 60.00000000000002
$ws.Range("A29").Value = -0.1976421426709416
$ws.Range("B29").Value = 59.99999999999999
$ws.Range("A30").Value = -0.185640942550929
$ws.Range("B30").Value = 60.00000000000002
$ws.Range("A31").Value = -0.1736397424309164
$ws.Range("B31").Value = 59.99999999999999
$ws.Range("A32").Value = -0.161638542310906
$ws.Range("B32").Value = 59.99999999999999
$ws.Range("A33").Value = -0.1496373421908936
$ws.Range("B33").Value = 59.99999999999999
$ws.Range("A34").Value = -0.1376361420708831
$ws.Range("B34").Value = 60.00000000000002
$ws.Range("A35").Value = -0.1256349419508704
$ws.Range("B35").Value = 59.99999999999999
$ws.Range("A36").Value = -0.1136337418308587
$ws.Range("B36").Value = 60.00000000000002
$ws.Range("A37").Value = -0.101632541710848
$ws.Range("B37").Value = 60.00000000000002
$ws.Range("A38").Value = -0.08963134159083529
$ws.Range("B38").Value = 60.00000000000002
$ws.Range("A39").Value = -0.07763014147082543
$ws.Range("B39").Value = 60.00000000000002
$ws.Range("A40").Value = -0.06562894135081271
$ws.Range("B40").Value = 59.99999999999999
$ws.Range("A41").Value = -0.05362774123080054
$ws.Range("B41").Value = 60.00000000000002
$ws.Range("A42").Value = -0.04162654111079032
$ws.Range("B42").Value = 60.00000000000002
$ws.Range("A43").Value = -0.02962534099077779
$ws.Range("B43").Value = 60.00000000000002
$ws.Range("A44").Value = -0.01762414087076759
$ws.Range("B44").Value = 60.00000000000002
$ws.Range("A45").Value = -0.005622940750755393
$ws.Range("B45").Value = 59.99999999999999
$ws.Range("A46").Value = 0.006378259369257228
$ws.Range("B46").Value = 59.99999999999999
$ws.Range("A47").Value = 0.01837945948926704
$ws.Range("B47").Value = 60.00000000000002
$ws.Range("A48").Value = 0.03038065960927907
$ws.Range("B48").Value = 60.00000000000002
$ws.Range("A49").Value = 0.04238185972928952
$ws.Range("B49").Value = 60.00000000000002
$ws.Range("A50").Value = 0.05438305984930228
$ws.Range("B50").Value = 59.99999999999999
$ws.Range("A51").Value = 0.06638425996931416
$ws.Range("B51").Value = 60.00000000000002
$ws.Range("A52").Value = 0.07838546008932472
$ws.Range("B52").Value = 59.99999999999999
$ws.Range("A53").Value = 0.09038666020933674
$ws.Range("B53").Value = 60.00000000000002
$ws.Range("A54").Value = 0.1023878603293474
$ws.Range("B54").Value = 59.99999999999999
$ws.Range("A55").Value = 0.1143890604493596
$ws.Range("B55").Value = 60.00000000000002
$ws.Range("A56").Value = 0.1263902605693717
$ws.Range("B56").Value = 60.00000000000002
$ws.Range("A57").Value = 0.1383914606893828
$ws.Range("B57").Value = 59.99999999999999
$ws.Range("A58").Value = 0.1503926608093951
$ws.Range("B58").Value = 60.00000000000002
$ws.Range("A59").Value = 0.162393860929406
$ws.Range("B59").Value = 59.99999999999999
$ws.Range("A60").Value = 0.1743950610494183
$ws.Range("B60").Value = 59.99999999999999
$ws.Range("A61").Value = 0.1863962611694308
$ws.Range("B61").Value = 59.99999999999999
$ws.Range("A62").Value = 0.1983974612894411
$ws.Range("B62").Value = 60.00000000000004
$ws.Range("A63").Value = 0.2103986614094541
$ws.Range("B63").Value = 59.99999999999999
$ws.Range("A64").Value = 0.2223998615294651
$ws.Range("B64").Value = 60.00000000000002
$ws.Range("A65").Value = 0.2344010616494781
$ws.Range("B65").Value = 59.99999999999999
$ws.Range("A66").Value = 0.2464022617695241
$ws.Range("B66").Value = 59.99999999999999
$ws.Range("A67").Value = 0.2584034618895207
$ws.Range("B67").Value = 59.99999999999999
$ws.Range("A68").Value = 0.2704046620095249
$ws.Range("B68").Value = 59.99999999999999
$ws.Range("A69").Value = 0.2824058621295311
$ws.Range("B69").Value = 59.99999999999999
$ws.Range("A70").Value = 0.2944070622495409
$ws.Range("B70").Value = 59.99999999999999
$ws.Range("A71").Value = 0.306408262369552
$ws.Range("B71").Value = 59.99999999999999
$ws.Range("A72").Value = 0.3184094624895625
$ws.Range("B72").Value = 60.00000000000002
$ws.Range("A73").Value = 0.3304106626095751
$ws.Range("B73").Value = 59.99999999999999
$ws.Range("A74").Value = 0.3424118627295863
$ws.Range("B74").Value = 59.99999999999999
$ws.Range("A75").Value = 0.3544130628495989
$ws.Range("B75").Value = 60.00000000000002
$ws.Range("A76").Value = 0.3664142629696119
$ws.Range("B76").Value = 59.99999999999999
$ws.Range("A77").Value = 0.3784154630896235
$ws.Range("B77").Value = 59.99999999999999
$ws.Range("A78").Value = 0.3904166632096365
$ws.Range("B78").Value = 59.99999999999999
$ws.Range("A79").Value = 0.4024178633296475
$ws.Range("B79").Value = 59.99999999999999
$ws.Range("A80").Value = 0.4144190634496606
$ws.Range("B80").Value = 60.00000000000002
$ws.Range("A81").Value = 0.4264202635696737
$ws.Range("B81").Value = 60.00000000000002
$ws.Range("A82").Value = 0.4384214636896847
$ws.Range("B82").Value = 59.99999999999999
$ws.Range("A83").Value = 0.4504226638096983
$ws.Range("B83").Value = 59.99999999999999
$ws.Range("A84").Value = 0.4624238639297171
$ws.Range("B84").Value = 60.00000000000002
$ws.Range("A85").Value = 0.4744250640497223
$ws.Range("B85").Value = 59.99999999999999
$ws.Range("A86").Value = 0.4864262641697362
$ws.Range("B86").Value = 59.99999999999999
$ws.Range("A87").Value = 0.4984274642897471
$ws.Range("B87").Value = 59.99999999999999
$ws.Range("A88").Value = 0.51042866440976
$ws.Range("B88").Value = 59.99999999999999
$ws.Range("A89").Value = 0.5224298645297716
$ws.Range("B89").Value = 60.00000000000002

# Remove now-unused trailing rows 90-92 to shrink dimension to A1:B89
$ws.Rows.Item(90).Delete()
$ws.Rows.Item(90).Delete()
$ws.Rows.Item(90).Delete()
